$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    above the existing 2022-Q2 row, shifting the rest down.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Copy formatting (styles) from the row that used to be row 2 (now row 3)
# down onto the newly inserted row 2, so the new row matches the look
# of the other data rows.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

# Fill in the values for the new 2022-Q3 row.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 12
$summary.Cells.Item(2, 4).Value = 1.55

# Renumber the running index in column A for the rows that shifted down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3

# ------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" sheet with fund holdings data,
#    positioned right before the existing "2022-Q2" sheet.
# ------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Duplicate the 2022-Q2 sheet (keeps header/column formatting identical)
# and place the copy immediately before it, then rename & refill it.
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item("2022-Q2 (2)")
$q3Sheet.Name = "2022-Q3"

# Clear out any leftover rows beyond what we need (the source sheet has
# 9 rows of data, the new sheet needs 13), then write fresh data.
$q3Sheet.Range("A2:H13").ClearContents()

# Rows 10-13 did not exist on the source sheet, so column A there is
# missing the bordered/bold index style used by the other rows. Copy
# that formatting down from row 2 onto the new rows.
$q3Sheet.Range("A2").Copy()
$q3Sheet.Range("A10:A13").PasteSpecial(-4122)

$data = @(
    @("004640", "华夏节能环保股票A", "6.04", "93.11", "8.53", "0.5152", 1),
    @("012703", "华夏核心成长混合A", "4.94", "93.34", "8.94", "0.4416", 1),
    @("015229", "华夏低碳经济一年持有混合A", "2.49", "91.97", "8.37", "0.2084", 1),
    @("014410", "华夏时代领航两年持有混合A", "2.65", "89.46", "4.22", "0.1118", 8),
    @("015230", "华夏低碳经济一年持有混合C", "1.02", "91.97", "8.37", "0.0854", 1),
    @("003300", "华夏圆和灵活配置混合A", "0.77", "75.31", "7.70", "0.0593", 1),
    @("012710", "华夏核心成长混合C", "0.61", "93.34", "8.94", "0.0545", 1),
    @("015068", "华夏圆和灵活配置混合C", "0.33", "75.31", "7.70", "0.0254", 1),
    @("015060", "华夏节能环保股票C", "0.26", "93.11", "8.53", "0.0222", 1),
    @("014411", "华夏时代领航两年持有混合C", "0.46", "89.46", "4.22", "0.0194", 8),
    @("002409", "华夏新活力灵活配置混合A", "0.15", "69.89", "3.10", "0.0046", 10),
    @("002410", "华夏新活力灵活配置混合C", "0.00", "69.89", "3.10", "", 10)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $q3Sheet.Cells.Item($row, 1).Value = $i
    $q3Sheet.Cells.Item($row, 2).Value = "'" + $item[0]
    $q3Sheet.Cells.Item($row, 3).Value = $item[1]
    $q3Sheet.Cells.Item($row, 4).Value = "'" + $item[2]
    $q3Sheet.Cells.Item($row, 5).Value = "'" + $item[3]
    $q3Sheet.Cells.Item($row, 6).Value = "'" + $item[4]
    $q3Sheet.Cells.Item($row, 8).Value = $item[6]
}

# Last row's "持有市值" column is stored as a real number (0) rather
# than the text "0", unlike the rows above it.
$q3Sheet.Cells.Item(13, 7).Value = 0
for ($i = 0; $i -lt ($data.Count - 1); $i++) {
    $row = $i + 2
    $item = $data[$i]
    $q3Sheet.Cells.Item($row, 7).Value = "'" + $item[5]
}

# Keep the originally-selected sheet ("2021-Q1", the last tab) active,
# rather than leaving the newly created sheet selected.
$wb.Worksheets.Item("2021-Q1").Activate()
